# "Fruta / hortaliza, semanal" — weekly refresh of the Haba (Feria Lagunitas
# de Puerto Montt) price series: a new week's observation is inserted at the
# top of the data block (row 23, right after the fixed first 21 data rows),
# pushing every later observation down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 23; Excel shifts rows 23:72 down to 24:73 and
# the sheet's dimension grows from R72 to R73 automatically.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row with this week's observation.
$ws.Range("A23").Value2 = 4
$ws.Range("B23").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C23").Value2 = "Los Lagos"
$ws.Range("D23").Value2 = 44544
$ws.Range("E23").Value2 = 10
$ws.Range("F23").Value2 = 100112026
$ws.Range("G23").Value2 = "Haba"
$ws.Range("H23").Value2 = "Sin especificar"
$ws.Range("I23").Value2 = "Primera"
$ws.Range("J23").Value2 = 120
$ws.Range("K23").Value2 = 12000
$ws.Range("L23").Value2 = 12000
$ws.Range("M23").Value2 = 12000
$ws.Range("N23").Value2 = "$/saco 25 kilos"
$ws.Range("O23").Value2 = "Región de La Araucanía"
$ws.Range("P23").Value2 = 480
$ws.Range("Q23").Value2 = 25
$ws.Range("R23").Value2 = "Hortaliza"
